$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 141, shifting existing rows 141:268 down to 142:269
$ws.Rows.Item(141).Insert()

# Populate the newly inserted row 141 with data
$ws.Cells.Item(141, 1).Value = 5
$ws.Cells.Item(141, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(141, 3).Value = "Maule"
$ws.Cells.Item(141, 4).Value = 44586
$ws.Cells.Item(141, 5).Value = 7
$ws.Cells.Item(141, 6).Value = 100114013
$ws.Cells.Item(141, 7).Value = "Zanahoria"
$ws.Cells.Item(141, 8).Value = "Sin especificar"
$ws.Cells.Item(141, 9).Value = "Primera"
$ws.Cells.Item(141, 10).Value = 300
$ws.Cells.Item(141, 11).Value = 8000
$ws.Cells.Item(141, 12).Value = 8000
$ws.Cells.Item(141, 13).Value = 8000
$ws.Cells.Item(141, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(141, 15).Value = "Región de Ñuble"
$ws.Cells.Item(141, 16).Value = 400
$ws.Cells.Item(141, 17).Value = 20
$ws.Cells.Item(141, 18).Value = "Hortaliza"
